# Update "想去人数" (F column) figures for the latest data refresh.
$wb = $excel.ActiveWorkbook

# 展览 (Exhibition) sheet
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 7892
$wsExpo.Range("F3").Value = 74
$wsExpo.Range("F5").Value = 59
$wsExpo.Range("F6").Value = 655
$wsExpo.Range("F7").Value = 1253
$wsExpo.Range("F8").Value = 213
$wsExpo.Range("F10").Value = 180

# 演出 (Performance) sheet
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 15

# 全部类型 (All types) sheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 7892
$wsAll.Range("F3").Value = 74
$wsAll.Range("F5").Value = 59
$wsAll.Range("F6").Value = 655
$wsAll.Range("F7").Value = 1253
$wsAll.Range("F8").Value = 213
$wsAll.Range("F9").Value = 15
$wsAll.Range("F11").Value = 180
